# Add NBs for direct data entry and recording query execution times.
#
# The benchmark sheet tracked query execution time (column C) against the
# number of records (column B) for three queries (Query 1/2/3, rows 2-4,
# 5-7, 8-10). The "large" record-count tier is changed from 10,000,000 to
# 5,000,000 rows, and the measured execution time for Query 2 at the new
# 5,000,000-row tier is recorded.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Direct data entry: update the large-size tier from 10,000,000 to 5,000,000
# rows for each of the three queries.
$ws.Range("B4").Value = 5000000
$ws.Range("B7").Value = 5000000
$ws.Range("B10").Value = 5000000

# Record the query execution time measured for Query 2 at the 5,000,000-row
# tier (Redis/MySQL/MongoDB timing column C).
$ws.Range("C7").Value = 8.0177612499999995

# Minimize the workbook's window, matching the authored workbook view state.
$wb.Windows.Item(1).WindowState = -4140
